$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New header cells for the two newly-logged columns (G = start_time,
#    H = end_time)
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "start_time"
$ws.Range("H1").Value = "end_time"

# ---------------------------------------------------------------------------
# 2. New data rows (32-37) with the pump-flow logging for 2020-11-01,
#    mirroring the existing A:E layout and adding the new G/H timestamps.
# ---------------------------------------------------------------------------
$rowDate = [DateTime]"2020-11-01"

# Row 32
$ws.Range("A32").Value = $rowDate
$ws.Range("B32").Value = "irr"
$ws.Range("C32").Value = "amb"
$ws.Range("D32").Value = 45
$ws.Range("E32").Value = 49
$ws.Range("G32").Value = 44501.031886574077
$ws.Range("H32").Value = 44501.03466435185

# Row 33
$ws.Range("A33").Value = $rowDate
$ws.Range("B33").Value = "irr"
$ws.Range("C33").Value = "cc"
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 7
$ws.Range("G33").Value = 44501
$ws.Range("H33").Value = 44501.005208333336

# Row 34
$ws.Range("A34").Value = $rowDate
$ws.Range("B34").Value = "irr"
$ws.Range("C34").Value = "oc"
$ws.Range("D34").Value = 8
$ws.Range("E34").Value = 11
$ws.Range("G34").Value = 44501.005902777775
$ws.Range("H34").Value = 44501.008275462962

# Row 35
$ws.Range("A35").Value = $rowDate
$ws.Range("B35").Value = "con"
$ws.Range("C35").Value = "cc"
$ws.Range("D35").Value = 12
$ws.Range("E35").Value = 19
$ws.Range("G35").Value = 44501.008333333331
$ws.Range("H35").Value = 44501.013541666667

# Row 36
$ws.Range("A36").Value = $rowDate
$ws.Range("B36").Value = "con"
$ws.Range("C36").Value = "oc"
$ws.Range("D36").Value = 20
$ws.Range("E36").Value = 27
$ws.Range("G36").Value = 44501.014236111114
$ws.Range("H36").Value = 44501.019386574073

# Row 37
$ws.Range("A37").Value = $rowDate
$ws.Range("B37").Value = "con"
$ws.Range("C37").Value = "amb"
$ws.Range("D37").Value = 45
$ws.Range("E37").Value = 49
$ws.Range("G37").Value = 44501.031886574077
$ws.Range("H37").Value = 44501.03466435185

# ---------------------------------------------------------------------------
# 3. Apply the date/time number format used for the new start_time/end_time
#    columns.
# ---------------------------------------------------------------------------
$ws.Range("G2:H37").NumberFormat = "m/d/yyyy\ h:mm:ss"

# ---------------------------------------------------------------------------
# 4. Resize the new columns to fit their content.
# ---------------------------------------------------------------------------
$ws.Columns("G:H").AutoFit()

# ---------------------------------------------------------------------------
# 5. Update the view: scroll so row 10 is at the top and leave the selection
#    on the next empty cell below the data (H38), matching where the user
#    left off after entering the new rows.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H38").Select()
